$d = $word.ActiveDocument

# The document contains four occurrences of an <id>...</id> tag whose
# three parts ("<id>", the identifier text, "</id>") were previously
# split across three separate runs. Collapse each into a single run
# (keeping the first run's formatting) containing the full tag text.
$ids = @("p104v_1", "p104v_2", "p104v_3", "p104v_4")

foreach ($id in $ids) {
    $old = "<id>" + $id + "</id>"
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $old, 2) | Out-Null
}
